$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "first"
$ws.Range("B6").Value = "last"
$ws.Range("C6").Value = "name"
$ws.Range("D6").Value = "pass"
